# initial refactor to OOP
# Populate the newly added Trade Activity summary columns (K:Q) for rows 4-13
# on the "Trade Activity" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trade Activity")

# Column layout:
#   K = Trade vs Ex-Ante %   (percentage, like F/G)
#   L = Ex-Post vs Ex-Ante % (percentage, like F/G)
#   M = Net Trade Value      (number, like C/D/E/H/I/J)
#   N = Net Buy Qty          (number)
#   O = Net Sell Qty         (number)
#   P = Net Buy Value        (number)
#   Q = Net Sell Value       (number)

$data = @{
    4  = @(0.1501646666616732, 1, 2998695.3, 31400, 0, 2998695.3, 0)
    5  = @(1, 1, 3377802036.354166, 33989000, 0, 3377802036.354166, 0)
    6  = @(0, 0, -15750, 0, 63, 0, 15750)
    7  = @(0, 0, -4697, 0, 7, 0, 4697)
    8  = @(0.07509845945615433, 0.992655645010849, 220671489.83, 2344178, 265928, 241802670.272, 21131180.442)
    9  = @(0, 0, -25783353, 0, 85048, 0, 25783353)
    10 = @(0.001542154835306817, 0.9984578451646932, 0, 0, 24000, 0, 0)
    11 = @(0, 0, -3696850, 0, 535, 0, 3696850)
    12 = @(0, 0, -985920.0000000001, 0, 312, 0, 985920.0000000001)
    13 = @(0, 0, -234940, 0, 34, 0, 234940)
}

foreach ($rowNum in $data.Keys) {
    $values = $data[$rowNum]

    $ws.Cells.Item($rowNum, 11).Value = $values[0]   # K
    $ws.Cells.Item($rowNum, 12).Value = $values[1]   # L
    $ws.Cells.Item($rowNum, 13).Value = $values[2]   # M
    $ws.Cells.Item($rowNum, 14).Value = $values[3]   # N
    $ws.Cells.Item($rowNum, 15).Value = $values[4]   # O
    $ws.Cells.Item($rowNum, 16).Value = $values[5]   # P
    $ws.Cells.Item($rowNum, 17).Value = $values[6]   # Q

    $ws.Range($ws.Cells.Item($rowNum, 11), $ws.Cells.Item($rowNum, 12)).NumberFormat = "0.00%"
    $ws.Range($ws.Cells.Item($rowNum, 13), $ws.Cells.Item($rowNum, 17)).NumberFormat = "#,##0.00"
}
